# Update gh-pages to output generated at 456a3b4
#
# This script updates the "想去人数" (F column) counters across all four
# worksheets, and marks the "上海·二次元小偶像之夜" event (展览!row18) as
# cancelled: its name gets a "（取消）" suffix and its "最低票价" (G column)
# switches from a numeric price to the text "不可售".

$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 ---
$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(3, 6).Value = 2411    # F3: 2409 -> 2411
$ws.Cells.Item(5, 6).Value = 139     # F5: 138 -> 139
$ws.Cells.Item(6, 6).Value = 63      # F6: 59 -> 63
$ws.Cells.Item(7, 6).Value = 274     # F7: 272 -> 274
$ws.Cells.Item(8, 6).Value = 334     # F8: 331 -> 334
$ws.Cells.Item(9, 6).Value = 2229    # F9: 2210 -> 2229
$ws.Cells.Item(10, 6).Value = 1159   # F10: 1156 -> 1159
$ws.Cells.Item(12, 6).Value = 851    # F12: 848 -> 851
$ws.Cells.Item(13, 6).Value = 91     # F13: 90 -> 91
$ws.Cells.Item(14, 6).Value = 842    # F14: 839 -> 842
$ws.Cells.Item(15, 6).Value = 1477   # F15: 1470 -> 1477
$ws.Cells.Item(16, 6).Value = 727    # F16: 718 -> 727
$ws.Cells.Item(17, 6).Value = 1696   # F17: 1691 -> 1696

# Row 18: "上海·二次元小偶像之夜" event cancelled
$ws.Cells.Item(18, 3).Value = "上海·二次元小偶像之夜（取消）"   # C18
$ws.Cells.Item(18, 7).Value = "不可售"                          # G18: 129 -> "不可售"

$ws.Cells.Item(19, 6).Value = 360    # F19: 355 -> 360
$ws.Cells.Item(20, 6).Value = 67     # F20: 65 -> 67
$ws.Cells.Item(21, 6).Value = 108    # F21: 104 -> 108
$ws.Cells.Item(23, 6).Value = 2624   # F23: 2615 -> 2624

# --- Sheet: 演出 ---
$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(10, 6).Value = 75     # F10: 74 -> 75
$ws.Cells.Item(19, 6).Value = 152    # F19: 151 -> 152
$ws.Cells.Item(22, 6).Value = 7      # F22: 6 -> 7
$ws.Cells.Item(26, 6).Value = 45     # F26: 44 -> 45
$ws.Cells.Item(28, 6).Value = 182    # F28: 180 -> 182
$ws.Cells.Item(31, 6).Value = 8      # F31: 7 -> 8
$ws.Cells.Item(38, 6).Value = 340    # F38: 339 -> 340
$ws.Cells.Item(46, 6).Value = 297    # F46: 296 -> 297

# --- Sheet: 本地生活 ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(4, 6).Value = 2484    # F4: 2482 -> 2484
$ws.Cells.Item(6, 6).Value = 2496    # F6: 2493 -> 2496
$ws.Cells.Item(8, 6).Value = 1843    # F8: 1842 -> 1843
$ws.Cells.Item(9, 6).Value = 134     # F9: 133 -> 134
$ws.Cells.Item(11, 6).Value = 15     # F11: 9 -> 15
$ws.Cells.Item(12, 6).Value = 357    # F12: 356 -> 357
$ws.Cells.Item(13, 6).Value = 2779   # F13: 2771 -> 2779
$ws.Cells.Item(15, 6).Value = 666    # F15: 664 -> 666

# --- Sheet: 全部类型 ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 2484    # F2: 2482 -> 2484
$ws.Cells.Item(4, 6).Value = 134     # F4: 133 -> 134
$ws.Cells.Item(7, 6).Value = 2411    # F7: 2409 -> 2411
$ws.Cells.Item(8, 6).Value = 2779    # F8: 2771 -> 2779
$ws.Cells.Item(10, 6).Value = 139    # F10: 138 -> 139
$ws.Cells.Item(11, 6).Value = 666    # F11: 664 -> 666
$ws.Cells.Item(16, 6).Value = 63     # F16: 59 -> 63
$ws.Cells.Item(17, 6).Value = 274    # F17: 272 -> 274
$ws.Cells.Item(18, 6).Value = 334    # F18: 331 -> 334
$ws.Cells.Item(21, 6).Value = 851    # F21: 848 -> 851
$ws.Cells.Item(22, 6).Value = 91     # F22: 90 -> 91
$ws.Cells.Item(23, 6).Value = 842    # F23: 839 -> 842
$ws.Cells.Item(28, 6).Value = 727    # F28: 718 -> 727
$ws.Cells.Item(30, 6).Value = 7      # F30: 6 -> 7
$ws.Cells.Item(31, 6).Value = 1696   # F31: 1691 -> 1696
$ws.Cells.Item(32, 6).Value = 360    # F32: 355 -> 360
$ws.Cells.Item(35, 6).Value = 45     # F35: 44 -> 45
$ws.Cells.Item(37, 6).Value = 182    # F37: 180 -> 182
$ws.Cells.Item(39, 6).Value = 67     # F39: 65 -> 67
$ws.Cells.Item(40, 6).Value = 8      # F40: 7 -> 8
$ws.Cells.Item(42, 6).Value = 340    # F42: 339 -> 340
$ws.Cells.Item(43, 6).Value = 2624   # F43: 2615 -> 2624
$ws.Cells.Item(48, 6).Value = 297    # F48: 296 -> 297
